$wb = $excel.ActiveWorkbook

# --- 1) ODI Batting: drop the empty placeholder cells in column B ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Range("B4").ClearContents()
$wsBatting.Range("B5").ClearContents()
$wsBatting.Range("B7").ClearContents()
$wsBatting.Range("B9").ClearContents()
$wsBatting.Range("B10").ClearContents()
$wsBatting.Range("B12").ClearContents()
$wsBatting.Range("B15").ClearContents()
$wsBatting.Range("B16").ClearContents()
$wsBatting.Range("B17").ClearContents()
$wsBatting.Range("B18").ClearContents()
$wsBatting.Range("B19").ClearContents()
$wsBatting.Range("B20").ClearContents()

# --- 2) Add the new "ODI Batting Extra" sheet as the last tab ---
$firstSheet = $wb.Worksheets.Item(1)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "ODI Batting Extra"

# Header row: reuse the bold/centered/bordered header style used by the other sheets
$firstSheet.Range("A1:D1").Copy()
$ws4.Range("A1:F1").PasteSpecial(-4122)
$ws4.Range("A1").Value = "MATCH_CODE"
$ws4.Range("B1").Value = "BATTING_POSITION"
$ws4.Range("C1").Value = "NUM_4"
$ws4.Range("D1").Value = "NUM_6"
$ws4.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws4.Range("F1").Value = "MAN_OF_MATCH"

# Data rows 2-20 (leading "'" forces text type, matching the scraped source data;
# bare "'" alone yields an explicit empty text cell rather than a blank one)
# Row 2
$ws4.Range("A2").Value = "'3188"
$ws4.Range("B2").Value = "'"
$ws4.Range("C2").Value = "'"
$ws4.Range("D2").Value = "'"
$ws4.Range("E2").Value = "'"
$ws4.Range("F2").Value = "'NO"
# Row 3
$ws4.Range("A3").Value = "'3190"
$ws4.Range("B3").Value = "'"
$ws4.Range("C3").Value = "'"
$ws4.Range("D3").Value = "'"
$ws4.Range("E3").Value = "'"
$ws4.Range("F3").Value = "'NO"
# Row 4
$ws4.Range("A4").Value = "'3220"
$ws4.Range("B4").Value = 11
$ws4.Range("C4").Value = "'"
$ws4.Range("D4").Value = "'"
$ws4.Range("E4").Value = "'"
$ws4.Range("F4").Value = "'NO"
# Row 5
$ws4.Range("A5").Value = "'3223"
$ws4.Range("B5").Value = 11
$ws4.Range("C5").Value = "'0"
$ws4.Range("D5").Value = "'0"
$ws4.Range("E5").Value = "'"
$ws4.Range("F5").Value = "'NO"
# Row 6
$ws4.Range("A6").Value = "'3225"
$ws4.Range("B6").Value = 11
$ws4.Range("C6").Value = "'0"
$ws4.Range("D6").Value = "'0"
$ws4.Range("E6").Value = "'0.96%"
$ws4.Range("F6").Value = "'NO"
# Row 7
$ws4.Range("A7").Value = "'3228"
$ws4.Range("B7").Value = "'"
$ws4.Range("C7").Value = "'"
$ws4.Range("D7").Value = "'"
$ws4.Range("E7").Value = "'"
$ws4.Range("F7").Value = "'NO"
# Row 8
$ws4.Range("A8").Value = "'3230"
$ws4.Range("B8").Value = "'"
$ws4.Range("C8").Value = "'"
$ws4.Range("D8").Value = "'"
$ws4.Range("E8").Value = "'"
$ws4.Range("F8").Value = "'NO"
# Row 9
$ws4.Range("A9").Value = "'3232"
$ws4.Range("B9").Value = "'"
$ws4.Range("C9").Value = "'"
$ws4.Range("D9").Value = "'"
$ws4.Range("E9").Value = "'"
$ws4.Range("F9").Value = "'NO"
# Row 10
$ws4.Range("A10").Value = "'3236"
$ws4.Range("B10").Value = "'"
$ws4.Range("C10").Value = "'"
$ws4.Range("D10").Value = "'"
$ws4.Range("E10").Value = "'"
$ws4.Range("F10").Value = "'"
# Row 11
$ws4.Range("A11").Value = "'3242"
$ws4.Range("B11").Value = "'"
$ws4.Range("C11").Value = "'"
$ws4.Range("D11").Value = "'"
$ws4.Range("E11").Value = "'"
$ws4.Range("F11").Value = "'"
# Row 12
$ws4.Range("A12").Value = "'3252"
$ws4.Range("B12").Value = "'"
$ws4.Range("C12").Value = "'"
$ws4.Range("D12").Value = "'"
$ws4.Range("E12").Value = "'"
$ws4.Range("F12").Value = "'"
# Row 13
$ws4.Range("A13").Value = "'3272"
$ws4.Range("B13").Value = "'"
$ws4.Range("C13").Value = "'"
$ws4.Range("D13").Value = "'"
$ws4.Range("E13").Value = "'"
$ws4.Range("F13").Value = "'"
# Row 14
$ws4.Range("A14").Value = "'3607"
$ws4.Range("B14").Value = "'"
$ws4.Range("C14").Value = "'"
$ws4.Range("D14").Value = "'"
$ws4.Range("E14").Value = "'"
$ws4.Range("F14").Value = "'"
# Row 15
$ws4.Range("A15").Value = "'3609"
$ws4.Range("B15").Value = "'"
$ws4.Range("C15").Value = "'"
$ws4.Range("D15").Value = "'"
$ws4.Range("E15").Value = "'"
$ws4.Range("F15").Value = "'"
# Row 16
$ws4.Range("A16").Value = "'4024"
$ws4.Range("B16").Value = "'"
$ws4.Range("C16").Value = "'"
$ws4.Range("D16").Value = "'"
$ws4.Range("E16").Value = "'"
$ws4.Range("F16").Value = "'"
# Row 17
$ws4.Range("A17").Value = "'4027"
$ws4.Range("B17").Value = "'"
$ws4.Range("C17").Value = "'"
$ws4.Range("D17").Value = "'"
$ws4.Range("E17").Value = "'"
$ws4.Range("F17").Value = "'"
# Row 18
$ws4.Range("A18").Value = "'4402"
$ws4.Range("B18").Value = "'"
$ws4.Range("C18").Value = "'"
$ws4.Range("D18").Value = "'"
$ws4.Range("E18").Value = "'"
$ws4.Range("F18").Value = "'"
# Row 19
$ws4.Range("A19").Value = "'4406"
$ws4.Range("B19").Value = "'"
$ws4.Range("C19").Value = "'"
$ws4.Range("D19").Value = "'"
$ws4.Range("E19").Value = "'"
$ws4.Range("F19").Value = "'"
# Row 20
$ws4.Range("A20").Value = "'4410"
$ws4.Range("B20").Value = "'"
$ws4.Range("C20").Value = "'"
$ws4.Range("D20").Value = "'"
$ws4.Range("E20").Value = "'"
$ws4.Range("F20").Value = "'"

$firstSheet.Activate()
Write-Output "edit applied: ODI Batting cleaned + ODI Batting Extra sheet added"
